$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 501; this shifts the existing rows 501-525 down to 502-526
$ws.Rows.Item(501).Insert()

# Populate the newly inserted row 501 with the new weekly price record
$ws.Range("A501").Value = 11
$ws.Range("B501").Value = "Vega Monumental Concepción"
$ws.Range("C501").Value = "Bíobío"
$ws.Range("D501").Value = 45147
$ws.Range("E501").Value = 8
$ws.Range("F501").Value = "Fruta"
$ws.Range("G501").Value = 100102
$ws.Range("H501").Value = "Cítricos"
$ws.Range("I501").Value = 100102005
$ws.Range("J501").Value = "Naranja"
$ws.Range("K501").Value = "Fukumoto"
$ws.Range("L501").Value = "Primera"
$ws.Range("M501").Value = 250
$ws.Range("N501").Value = 6000
$ws.Range("O501").Value = 6500
$ws.Range("P501").Value = 6200
$ws.Range("Q501").Value = "`$/bandeja 15 kilos empedrada"
$ws.Range("R501").Value = "Región de O'Higgins"
$ws.Range("S501").Value = 413
$ws.Range("T501").Value = 15
